$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3766981.8
$ws.Range("I43").Value = 8778629
$ws.Range("J43").Value = 8246.375
$ws.Range("K43").Value = 8778629
$ws.Range("L43").Value = 8246.375
$ws.Range("M43").Value = -8778560
$ws.Range("N43").Value = -8384.375
$ws.Range("H82").Value = 2600
$ws.Range("I82").Value = 2600
$ws.Range("K82").Value = 7800
$ws.Range("M82").Value = -7394
$ws.Range("H85").Value = 2600
$ws.Range("I85").Value = 2600
$ws.Range("K85").Value = 7800
$ws.Range("M85").Value = -6396
$ws.Range("H86").Value = 4000.6667
$ws.Range("I86").Value = 999
$ws.Range("J86").Value = 5501.5
$ws.Range("K86").Value = 999
$ws.Range("L86").Value = 5501.5
$ws.Range("M86").Value = 124
$ws.Range("N86").Value = -7747.5
$ws.Range("H89").Value = 4000.6667
$ws.Range("I89").Value = 999
$ws.Range("J89").Value = 5501.5
$ws.Range("K89").Value = 4995
$ws.Range("L89").Value = 27507.5
$ws.Range("M89").Value = 621
$ws.Range("N89").Value = -38739.5
$ws.Range("H96").Value = 61.555557
$ws.Range("I96").Value = 58.166668
$ws.Range("J96").Value = 68.333336
$ws.Range("K96").Value = 174.500004
$ws.Range("L96").Value = 205.000008
$ws.Range("M96").Value = 1198.499996
$ws.Range("N96").Value = -2951.000008
$ws.Range("H132").Value = 30308814
$ws.Range("I132").Value = 37042660
$ws.Range("K132").Value = 111127980
$ws.Range("M132").Value = -111125450
$ws.Range("H141").Value = 2566.5
$ws.Range("J141").Value = 2749.2222
$ws.Range("L141").Value = 8247.6666
$ws.Range("N141").Value = -18607.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1475.3066
$ws.Range("I32").Value = 1422.1803
$ws.Range("J32").Value = 1706.7858
$ws.Range("K32").Value = 1422.1803
$ws.Range("L32").Value = 1706.7858
$ws.Range("M32").Value = -1135.1803
$ws.Range("N32").Value = -2280.7858
$ws.Range("H45").Value = 1537.7858
$ws.Range("J45").Value = 1861
$ws.Range("L45").Value = 1861
$ws.Range("N45").Value = -2615
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()
$ws.Range("H61").Value = 6651.5454
$ws.Range("I61").Value = 4964.9473
$ws.Range("K61").Value = 4964.9473
$ws.Range("M61").Value = -4752.9473
$ws.Range("H110").Value = 4859.6
$ws.Range("I110").Value = 2949.75
$ws.Range("K110").Value = 2949.75
$ws.Range("M110").Value = -904.75
$ws.Range("H122").Value = 2757.3635
$ws.Range("I122").Value = 2720.2222
$ws.Range("J122").Value = 2924.5
$ws.Range("K122").Value = 8160.6666
$ws.Range("L122").Value = 8773.5
$ws.Range("M122").Value = -5710.6666
$ws.Range("N122").Value = -13673.5
$ws.Range("H136").Value = 6651.5454
$ws.Range("I136").Value = 4964.9473
$ws.Range("K136").Value = 14894.8419
$ws.Range("M136").Value = -12344.8419

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 816.26086
$ws.Range("I20").Value = 806.5789
$ws.Range("K20").Value = 806.5789
$ws.Range("M20").Value = -559.5789
$ws.Range("H86").Value = 3532485.2
$ws.Range("I86").Value = 5151634
$ws.Range("J86").Value = 24329.834
$ws.Range("K86").Value = 5151634
$ws.Range("L86").Value = 24329.834
$ws.Range("M86").Value = -5150511
$ws.Range("N86").Value = -26575.834
$ws.Range("H89").Value = 3532485.2
$ws.Range("I89").Value = 5151634
$ws.Range("J89").Value = 24329.834
$ws.Range("K89").Value = 25758170
$ws.Range("L89").Value = 121649.17
$ws.Range("M89").Value = -25752554
$ws.Range("N89").Value = -132881.17
$ws.Range("H94").Value = 2633.4
$ws.Range("I94").Value = 2750.4443
$ws.Range("K94").Value = 2750.4443
$ws.Range("M94").Value = -2299.4443
$ws.Range("H107").Value = 2315.45
$ws.Range("I107").Value = 2384.6843
$ws.Range("K107").Value = 2384.6843
$ws.Range("M107").Value = -464.6842999999999
$ws.Range("H134").Value = 3209.5652
$ws.Range("I134").Value = 3209.5652
$ws.Range("K134").Value = 9628.695599999999
$ws.Range("M134").Value = -7093.695599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1700.4
$ws.Range("I16").Value = 1697
$ws.Range("J16").Value = 1705.5
$ws.Range("K16").Value = 1697
$ws.Range("L16").Value = 1705.5
$ws.Range("M16").Value = -1410
$ws.Range("N16").Value = -2279.5
$ws.Range("H50").Value = 35125
$ws.Range("I50").Value = 250
$ws.Range("J50").Value = 70000
$ws.Range("K50").Value = 250
$ws.Range("L50").Value = 70000
$ws.Range("M50").Value = 375
$ws.Range("N50").Value = -71250
$ws.Range("H51").Value = 55999.4
$ws.Range("I51").Value = 46666
$ws.Range("J51").Value = 69999.5
$ws.Range("K51").Value = 46666
$ws.Range("L51").Value = 69999.5
$ws.Range("M51").Value = -45930
$ws.Range("N51").Value = -71471.5
$ws.Range("H60").Value = 69666
$ws.Range("I60").Value = 68999
$ws.Range("J60").Value = 69999.5
$ws.Range("K60").Value = 68999
$ws.Range("L60").Value = 69999.5
$ws.Range("M60").Value = -68488
$ws.Range("N60").Value = -71021.5
$ws.Range("H61").Value = 55999.4
$ws.Range("I61").Value = 46666
$ws.Range("J61").Value = 69999.5
$ws.Range("K61").Value = 46666
$ws.Range("L61").Value = 69999.5
$ws.Range("M61").Value = -46318
$ws.Range("N61").Value = -70695.5
$ws.Range("H68").Value = 56196.668
$ws.Range("I68").Value = 48000
$ws.Range("J68").Value = 60295
$ws.Range("K68").Value = 48000
$ws.Range("L68").Value = 60295
$ws.Range("M68").Value = -47251
$ws.Range("N68").Value = -61793
$ws.Range("H71").Value = 56196.668
$ws.Range("I71").Value = 48000
$ws.Range("J71").Value = 60295
$ws.Range("K71").Value = 144000
$ws.Range("L71").Value = 180885
$ws.Range("M71").Value = -140256
$ws.Range("N71").Value = -188373
$ws.Range("H99").Value = 5162.778
$ws.Range("I99").Value = 4640.467
$ws.Range("K99").Value = 4640.467
$ws.Range("M99").Value = -3142.467
$ws.Range("H109").Value = 50918.668
$ws.Range("J109").Value = 50918.668
$ws.Range("L109").Value = 50918.668
$ws.Range("N109").Value = -52998.668
$ws.Range("H113").Value = 1700.4
$ws.Range("I113").Value = 1697
$ws.Range("J113").Value = 1705.5
$ws.Range("K113").Value = 1697
$ws.Range("L113").Value = 1705.5
$ws.Range("M113").Value = 473
$ws.Range("N113").Value = -6045.5
$ws.Range("H126").Value = 5162.778
$ws.Range("I126").Value = 4640.467
$ws.Range("K126").Value = 13921.401
$ws.Range("M126").Value = -11451.401
$ws.Range("H132").Value = 1441
$ws.Range("I132").Value = 1441
$ws.Range("K132").Value = 4323
$ws.Range("M132").Value = -1793
$ws.Range("H134").Value = 1681.8334
$ws.Range("I134").Value = 1681.8334
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5045.5002
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2510.5002
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 2416.9
$ws.Range("I18").Value = 363.33334
$ws.Range("K18").Value = 1090.00002
$ws.Range("M18").Value = -921.0000199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 90.375
$ws.Range("I9").Value = 87.5
$ws.Range("J9").Value = 99
$ws.Range("K9").Value = 87.5
$ws.Range("L9").Value = 99
$ws.Range("M9").Value = 82.5
$ws.Range("N9").Value = -439
$ws.Range("H70").Value = 8604.464
$ws.Range("I70").Value = 6809.5293
$ws.Range("K70").Value = 6809.5293
$ws.Range("M70").Value = -6539.5293
$ws.Range("H73").Value = 8604.464
$ws.Range("I73").Value = 6809.5293
$ws.Range("K73").Value = 6809.5293
$ws.Range("M73").Value = -5873.5293
$ws.Range("H80").Value = 4335.3335
$ws.Range("J80").Value = 5339.2144
$ws.Range("L80").Value = 5339.2144
$ws.Range("N80").Value = -7335.2144
$ws.Range("H83").Value = 4335.3335
$ws.Range("J83").Value = 5339.2144
$ws.Range("L83").Value = 26696.072
$ws.Range("N83").Value = -36680.072
$ws.Range("H122").Value = 200206000
$ws.Range("I122").Value = 250254990
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 750764970
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -750762520
$ws.Range("N122").Value = -34900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 4950
$ws.Range("J13").Value = 4950
$ws.Range("L13").Value = 4950
$ws.Range("N13").Value = -5230
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H70").Value = 26185
$ws.Range("I70").Value = 30148
$ws.Range("K70").Value = 30148
$ws.Range("M70").Value = -29878
$ws.Range("H73").Value = 26185
$ws.Range("I73").Value = 30148
$ws.Range("K73").Value = 30148
$ws.Range("M73").Value = -29212
$ws.Range("H93").Value = 836.73334
$ws.Range("I93").Value = 545.5454999999999
$ws.Range("K93").Value = 545.5454999999999
$ws.Range("M93").Value = 702.4545000000001
$ws.Range("H132").Value = 9825.666999999999
$ws.Range("I132").Value = 9812.5
$ws.Range("J132").Value = 9852
$ws.Range("K132").Value = 29437.5
$ws.Range("L132").Value = 29556
$ws.Range("M132").Value = -26907.5
$ws.Range("N132").Value = -34616

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7086
$ws.Range("J81").Value = 8492.143
$ws.Range("L81").Value = 16984.286
$ws.Range("N81").Value = -19106.286
$ws.Range("H84").Value = 7086
$ws.Range("J84").Value = 8492.143
$ws.Range("L84").Value = 84921.42999999999
$ws.Range("N84").Value = -95529.42999999999
$ws.Range("H132").Value = 2067.1667
$ws.Range("I132").Value = 1980.6
$ws.Range("K132").Value = 5941.799999999999
$ws.Range("M132").Value = -3411.799999999999
